$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.28031864006898
$ws.Range("C2").Value = 10.92053664895577
$ws.Range("D2").Value = 5.060179375518758
$ws.Range("E2").Value = 12.55291291216699
$ws.Range("F2").Value = 25.98620036343339
$ws.Range("L2").Value = 9.860504489982578
$ws.Range("M2").Value = 14.82307386167019
$ws.Range("N2").Value = 18.24257211158297
$ws.Range("O2").Value = 23.16192675530845

$ws.Range("B3").Value = 14.84093108786048
$ws.Range("C3").Value = 10.75685622766351
$ws.Range("D3").Value = 5.032157787610223
$ws.Range("E3").Value = 12.59217483257685
$ws.Range("F3").Value = 25.95894039593519
$ws.Range("L3").Value = 9.86796940765001
$ws.Range("M3").Value = 14.737808707223
$ws.Range("N3").Value = 18.30153048265195
$ws.Range("O3").Value = 23.19811678039142

$ws.Range("B4").Value = 14.56694563306105
$ws.Range("C4").Value = 10.65417043017142
$ws.Range("D4").Value = 5.014649981082597
$ws.Range("E4").Value = 12.61762104840495
$ws.Range("F4").Value = 25.94989977007697
$ws.Range("L4").Value = 9.873900732569469
$ws.Range("M4").Value = 14.68755931435057
$ws.Range("N4").Value = 18.33958109138847
$ws.Range("O4").Value = 23.22626614977709

$ws.Range("B5").Value = 14.45441237404889
$ws.Range("C5").Value = 10.61180508330514
$ws.Range("D5").Value = 5.007441558837296
$ws.Range("E5").Value = 12.62832814387292
$ws.Range("F5").Value = 25.94815357830808
$ws.Range("L5").Value = 9.876657197928937
$ws.Range("M5").Value = 14.66762743703651
$ws.Range("N5").Value = 18.35555340575683
$ws.Range("O5").Value = 23.23922446858853

$ws.Range("B6").Value = 14.43567834406985
$ws.Range("C6").Value = 10.6047398180716
$ws.Range("D6").Value = 5.006240209850445
$ws.Range("E6").Value = 12.63012645987443
$ws.Range("F6").Value = 25.94798070707382
$ws.Range("L6").Value = 9.877135418629786
$ws.Range("M6").Value = 14.66435114725347
$ws.Range("N6").Value = 18.35823380028428
$ws.Range("O6").Value = 23.24146589712077

$ws.Range("B7").Value = 14.56543130318199
$ws.Range("C7").Value = 10.65360114238102
$ws.Range("D7").Value = 5.014553061891352
$ws.Range("E7").Value = 12.61776408022724
$ws.Range("F7").Value = 25.94986837183609
$ws.Range("L7").Value = 9.873936532414616
$ws.Range("M7").Value = 14.68728827799356
$ws.Range("N7").Value = 18.33979460964699
$ws.Range("O7").Value = 23.22643489363529

$ws.Range("B8").Value = 15.12979048646461
$ws.Range("C8").Value = 10.86457061348768
$ws.Range("D8").Value = 5.050581620160253
$ws.Range("E8").Value = 12.56617294193402
$ws.Range("F8").Value = 25.97520559019668
$ws.Range("L8").Value = 9.862798955286911
$ws.Range("M8").Value = 14.79324771467093
$ws.Range("N8").Value = 18.26251772875186
$ws.Range("O8").Value = 23.17317262950706

$ws.Range("B9").Value = 16.19609998333182
$ws.Range("C9").Value = 11.25969596024865
$ws.Range("D9").Value = 5.118738642219551
$ws.Range("E9").Value = 12.47559282237686
$ws.Range("F9").Value = 26.0857779604875
$ws.Range("L9").Value = 9.851629899209149
$ws.Range("M9").Value = 15.01698409123924
$ws.Range("N9").Value = 18.125601281507
$ws.Range("O9").Value = 23.11590196647006

$ws.Range("B10").Value = 16.94625651222384
$ws.Range("C10").Value = 11.53707234610248
$ws.Range("D10").Value = 5.167167530514305
$ws.Range("E10").Value = 12.41544858362888
$ws.Range("F10").Value = 26.20376620621698
$ws.Range("L10").Value = 9.849896505448648
$ws.Range("M10").Value = 15.19005049172731
$ws.Range("N10").Value = 18.03384609527445
$ws.Range("O10").Value = 23.10273681944911

$ws.Range("B11").Value = 17.27871705337084
$ws.Range("C11").Value = 11.66013701608393
$ws.Range("D11").Value = 5.188815307826541
$ws.Range("E11").Value = 12.38946742029878
$ws.Range("F11").Value = 26.26530547211666
$ws.Range("L11").Value = 9.850504878624486
$ws.Range("M11").Value = 15.27043723144002
$ws.Range("N11").Value = 17.99400693772245
$ws.Range("O11").Value = 23.10304496856792

$ws.Range("B12").Value = 17.40323155515513
$ws.Range("C12").Value = 11.70626557082461
$ws.Range("D12").Value = 5.196955464099913
$ws.Range("E12").Value = 12.37982649614589
$ws.Range("F12").Value = 26.28972708571107
$ws.Range("L12").Value = 9.850935257959692
$ws.Range("M12").Value = 15.30109579543419
$ws.Range("N12").Value = 17.97919302236469
$ws.Range("O12").Value = 23.10406758069074

$ws.Range("B13").Value = 17.37647829724856
$ws.Range("C13").Value = 11.69635240514009
$ws.Range("D13").Value = 5.195204925244421
$ws.Range("E13").Value = 12.38189406510102
$ws.Range("F13").Value = 26.28441795668706
$ws.Range("L13").Value = 9.850833687664911
$ws.Range("M13").Value = 15.2944835668582
$ws.Range("N13").Value = 17.98237136901768
$ws.Range("O13").Value = 23.10380705254732

$ws.Range("B14").Value = 17.28898924573603
$ws.Range("C14").Value = 11.66394166794676
$ws.Range("D14").Value = 5.189486166872435
$ws.Range("E14").Value = 12.38867029978437
$ws.Range("F14").Value = 26.26729232720161
$ws.Range("L14").Value = 9.850536283134142
$ws.Range("M14").Value = 15.27295528279235
$ws.Range("N14").Value = 17.99278273625362
$ws.Range("O14").Value = 23.10311094355096

$ws.Range("B15").Value = 17.23521653595523
$ws.Range("C15").Value = 11.64402677736181
$ws.Range("D15").Value = 5.185975719458949
$ws.Range("E15").Value = 12.39284664959589
$ws.Range("F15").Value = 26.25694756666055
$ws.Range("L15").Value = 9.850380133170122
$ws.Range("M15").Value = 15.25979634208009
$ws.Range("N15").Value = 17.9991954322046
$ws.Range("O15").Value = 23.10280253496498

$ws.Range("B16").Value = 16.92434228651574
$ws.Range("C16").Value = 11.5289648940549
$ws.Range("D16").Value = 5.165744894143677
$ws.Range("E16").Value = 12.41717419906017
$ws.Range("F16").Value = 26.19990167506954
$ws.Range("L16").Value = 9.849884783420897
$ws.Range("M16").Value = 15.18482864140866
$ws.Range("N16").Value = 18.03648783840173
$ws.Range("O16").Value = 23.10284344524714

$ws.Range("B17").Value = 16.73129274875405
$ws.Range("C17").Value = 11.45756179810276
$ws.Range("D17").Value = 5.153234227316891
$ws.Range("E17").Value = 12.43245100431858
$ws.Range("F17").Value = 26.16691177085898
$ws.Range("L17").Value = 9.849938087996136
$ws.Range("M17").Value = 15.13924850643915
$ws.Range("N17").Value = 18.05985161768624
$ws.Range("O17").Value = 23.10448188293381

$ws.Range("B18").Value = 16.6194347103233
$ws.Range("C18").Value = 11.41620084453282
$ws.Range("D18").Value = 5.146002534002359
$ws.Range("E18").Value = 12.44136763368883
$ws.Range("F18").Value = 26.14867811681939
$ws.Range("L18").Value = 9.850100287624082
$ws.Range("M18").Value = 15.11318932803747
$ws.Range("N18").Value = 18.07346881128057
$ws.Range("O18").Value = 23.10601699523926

$ws.Range("B19").Value = 16.58142424528968
$ws.Range("C19").Value = 11.40214741711114
$ws.Range("D19").Value = 5.143547914787285
$ws.Range("E19").Value = 12.44440896805152
$ws.Range("F19").Value = 26.14263219666307
$ws.Range("L19").Value = 9.850177823040754
$ws.Range("M19").Value = 15.1043937771994
$ws.Range("N19").Value = 18.07811012784884
$ws.Range("O19").Value = 23.10663853512571

$ws.Range("B20").Value = 16.75192906277283
$ws.Range("C20").Value = 11.46519316581709
$ws.Range("D20").Value = 5.154569740545899
$ws.Range("E20").Value = 12.43081133193874
$ws.Range("F20").Value = 26.17034697311213
$ws.Range("L20").Value = 9.84991880577142
$ws.Range("M20").Value = 15.14408445657686
$ws.Range("N20").Value = 18.05734598862023
$ws.Range("O20").Value = 23.10424611824225

$ws.Range("B21").Value = 17.31472526009174
$ws.Range("C21").Value = 11.67347453021112
$ws.Range("D21").Value = 5.191167480896461
$ws.Range("E21").Value = 12.38667459858298
$ws.Range("F21").Value = 26.27229230972749
$ws.Range("L21").Value = 9.850618217406172
$ws.Range("M21").Value = 15.27927291469536
$ws.Range("N21").Value = 17.98971728291943
$ws.Range("O21").Value = 23.10329082131135

$ws.Range("B22").Value = 17.67444871893892
$ws.Range("C22").Value = 11.80682842769713
$ws.Range("D22").Value = 5.214750650228784
$ws.Range("E22").Value = 12.35898001991687
$ws.Range("F22").Value = 26.34542889778534
$ws.Range("L22").Value = 9.852240625375268
$ws.Range("M22").Value = 15.3688864464927
$ws.Range("N22").Value = 17.94710478287547
$ws.Range("O22").Value = 23.10794658670709

$ws.Range("B23").Value = 17.48323412127834
$ws.Range("C23").Value = 11.73591646108008
$ws.Range("D23").Value = 5.202195343185988
$ws.Range("E23").Value = 12.37365601483699
$ws.Range("F23").Value = 26.30580368864528
$ws.Range("L23").Value = 9.851268398918332
$ws.Range("M23").Value = 15.32094959183528
$ws.Range("N23").Value = 17.96970301671952
$ws.Range("O23").Value = 23.10497863107964

$ws.Range("B24").Value = 16.74260209002846
$ws.Range("C24").Value = 11.46174398927726
$ws.Range("D24").Value = 5.153966076754709
$ws.Range("E24").Value = 12.43155221102696
$ws.Range("F24").Value = 26.16879163492854
$ws.Range("L24").Value = 9.849927113423993
$ws.Range("M24").Value = 15.14189766888968
$ws.Range("N24").Value = 18.05847820701938
$ws.Range("O24").Value = 23.10435085983831

$ws.Range("B25").Value = 15.91289847351181
$ws.Range("C25").Value = 11.15496236602506
$ws.Range("D25").Value = 5.100580398535953
$ws.Range("E25").Value = 12.49896871819288
$ws.Range("F25").Value = 26.04937746147631
$ws.Range("L25").Value = 9.853512081559842
$ws.Range("M25").Value = 14.95485454697831
$ws.Range("N25").Value = 18.16108336179724
$ws.Range("O25").Value = 23.12632647971856
